$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Typography sheet: add a new font entry "IP_Address" (row 7)
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

# Reset style on the cells we are about to populate so they don't inherit
# the column default style (mirrors how the source file has no explicit
# cell-level style override on these rows).
$typo.Range("B7:H7").Style = "Normal"

$typo.Cells.Item(7, 2).Value = "IP_Address"   # Typography Name
$typo.Cells.Item(7, 3).Value = "arial.ttf"    # Font
$typo.Cells.Item(7, 4).Value = 15             # Size
$typo.Cells.Item(7, 5).Value = 4              # Bpp
$typo.Cells.Item(7, 6).Value = "+"            # Fallback Character
$typo.Cells.Item(7, 7).Value = "."            # Wildcard Characters
$typo.Cells.Item(7, 8).Value = "0-9"          # Wildcard Ranges

# Touch (without changing) column I on rows 7 & 8 so an empty placeholder
# cell is materialised there too, matching the table's auto-fill behaviour.
$typo.Cells.Item(7, 9).Font.Bold = $false
$typo.Cells.Item(8, 9).Font.Bold = $false

# ---------------------------------------------------------------------------
# Translation sheet: new GB-ALIGNMENT / GB-DIRECTION columns, remove two
# one-off rows, and append the new IP-address related text rows
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# New header cells for the per-language alignment/direction override columns
$tr.Range("G3:H3").Style = "Normal"
$tr.Cells.Item(3, 7).Value = "GB-ALIGNMENT"
$tr.Cells.Item(3, 8).Value = "GB-DIRECTION"

# Rows 5-18 shift up by one: the old "SingleUseId1" row (old row 5) and the
# old "SingleUseId21" / 192.168.0.17 row (old row 18) are dropped, so every
# row from the old row 6 through old row 19 moves up one position.
$tr.Cells.Item(5, 2).Value = "SingleUseId3"
$tr.Cells.Item(5, 3).Value = "ButtonDown"
$tr.Cells.Item(5, 4).Value = "Center"
$tr.Cells.Item(5, 5).Value = "Reset"

$tr.Cells.Item(6, 2).Value = "toggleButtonWhite"
$tr.Cells.Item(6, 3).Value = "Label"
$tr.Cells.Item(6, 4).Value = "Center"
$tr.Cells.Item(6, 5).Value = "White Text"

$tr.Cells.Item(7, 2).Value = "toggleButtonOrange"
$tr.Cells.Item(7, 3).Value = "Label"
$tr.Cells.Item(7, 4).Value = "Center"
$tr.Cells.Item(7, 5).Value = "Orange Text"

$tr.Cells.Item(8, 2).Value = "SingleUseId10"
$tr.Cells.Item(8, 3).Value = "Label"
$tr.Cells.Item(8, 4).Value = "Center"
$tr.Cells.Item(8, 5).Value = "Send TCP"

$tr.Cells.Item(9, 2).Value = "SingleUseId13"
$tr.Cells.Item(9, 3).Value = "Label"
$tr.Cells.Item(9, 4).Value = "Center"
$tr.Cells.Item(9, 5).Value = "Toggle LED"

$tr.Cells.Item(10, 2).Value = "SingleUseId14"
$tr.Cells.Item(10, 3).Value = "Label"
$tr.Cells.Item(10, 4).Value = "Center"
$tr.Cells.Item(10, 5).Value = "Log Data"

$tr.Cells.Item(11, 2).Value = "SingleUseId15"
$tr.Cells.Item(11, 3).Value = "Label"
$tr.Cells.Item(11, 4).Value = "Center"
$tr.Cells.Item(11, 5).Value = "Dump Log"

$tr.Cells.Item(12, 2).Value = "SingleUseId16"
$tr.Cells.Item(12, 3).Value = "Label"
$tr.Cells.Item(12, 4).Value = "Center"
$tr.Cells.Item(12, 5).Value = "Clear Log"

$tr.Cells.Item(13, 2).Value = "SingleUseId17"
$tr.Cells.Item(13, 3).Value = "Default"
$tr.Cells.Item(13, 4).Value = "Left"
$tr.Cells.Item(13, 5).Value = "25"

$tr.Cells.Item(14, 2).Value = "SingleUseId18"
$tr.Cells.Item(14, 3).Value = "Label"
$tr.Cells.Item(14, 4).Value = "Center"
$tr.Cells.Item(14, 5).Value = "Random"

$tr.Cells.Item(15, 2).Value = "SingleUseId19"
$tr.Cells.Item(15, 3).Value = "Label"
$tr.Cells.Item(15, 4).Value = "Center"
$tr.Cells.Item(15, 5).Value = "Send TCP"

$tr.Cells.Item(16, 2).Value = "SingleUseId20"
$tr.Cells.Item(16, 3).Value = "Label"
$tr.Cells.Item(16, 4).Value = "Right"
$tr.Cells.Item(16, 5).Value = "IP Address"

$tr.Cells.Item(17, 2).Value = "SingleUseId22"
$tr.Cells.Item(17, 3).Value = "Label"
$tr.Cells.Item(17, 4).Value = "Center"
$tr.Cells.Item(17, 5).Value = "Set Zero"

$tr.Cells.Item(18, 2).Value = "SingleUseId23"
$tr.Cells.Item(18, 3).Value = "Default"
$tr.Cells.Item(18, 4).Value = "Left"
$tr.Cells.Item(18, 5).Value = "25"

# Row 19: becomes the new "ipAddrText" / IP_Address entry, plus it gets the
# two new per-language override cells (G19/H19)
$tr.Range("G19:H19").Style = "Normal"
$tr.Cells.Item(19, 2).Value = "ipAddrText"
$tr.Cells.Item(19, 3).Value = "IP_Address"
$tr.Cells.Item(19, 4).Value = "Left"
$tr.Cells.Item(19, 5).Value = "<ipaddress>"
$tr.Cells.Item(19, 7).Value = "Left"
$tr.Cells.Item(19, 8).Value = "LTR"

# Row 20: existing row, only the Text ID changes
$tr.Cells.Item(20, 2).Value = "SingleUseId24"

# Row 21: brand new row
$tr.Range("B21:F21").Style = "Normal"
$tr.Cells.Item(21, 2).Value = "SingleUseId25"
$tr.Cells.Item(21, 3).Value = "IP_Address"
$tr.Cells.Item(21, 4).Value = "Left"
$tr.Cells.Item(21, 5).Value = "0.0.0.0"
$tr.Cells.Item(21, 6).Value = "LTR"
